# Adapt column header formatting to respective input file names:
#   "<field>_old" -> "<field>_FV2310"
#   "<field>_new" -> "<field>_FV2404"
# Then (re)build the "Table1" structured table over the header so the new
# names become the table's column headers, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fieldNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J (1-10): "<field>_old" -> "<field>_FV2310"
for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($fieldNames[$i])_FV2310"
}

# Column K (11) is "diff" and is left untouched.

# Columns L-U (12-21): "<field>_new" -> "<field>_FV2404"
for ($i = 0; $i -lt $fieldNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($fieldNames[$i])_FV2404"
}

# Turn the header + data range into a structured table ("Table1") so the
# renamed headers also become the table's column names.
$tableRange = $ws.Range("A1:U68")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split below row 1).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
